# TC_5-FDR_E2E-FDR-2938-SYSTEM ADMIN REFERENCE DATA
#
# Update the voucher's draw date (B2). The paydate (C2, =B2+2) and the
# TODAY()-driven reference-date columns (N2:Q2, AD2) recalculate
# automatically from the pinned clock.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44116

# Re-point the window at the System Admin Reference Data columns (N:Q)
# and leave Q19 as the active cell, matching where the reviewer left off.
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 1
$ws.Range("Q19").Select()
